$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, pushing existing rows 70-86 down to 71-87.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly data entry.
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value = "Bíobío"
$ws.Cells.Item(70, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = 100112012
$ws.Cells.Item(70, 7).Value = "Espinaca"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 150
$ws.Cells.Item(70, 11).Value = 6500
$ws.Cells.Item(70, 12).Value = 7000
$ws.Cells.Item(70, 13).Value = 6733
$ws.Cells.Item(70, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 673
$ws.Cells.Item(70, 17).Value = 10
$ws.Cells.Item(70, 18).Value = "Hortaliza"
